# Update scripts with new TPM values.
# The old "ECs -> *" rows (2:4) are dropped; the former "FAPs -> *" rows
# (5:7) move up to become rows 2:4, and several of their metric columns
# are refreshed with newly-computed TPM-derived numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the obsolete "ECs" sending-cluster rows; rows 5:7 ("FAPs" sending)
# shift up to rows 2:4.
$ws.Rows("2:4").Delete()

# Refresh the recomputed metrics for the new row 2 (FAPs -> ECs)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.07259900000000001
$ws.Range("N2").Value = 0.217797
$ws.Range("O2").Value = 0.0162094769588191
$ws.Range("P2").Value = 0.0162094769588191
$ws.Range("Q2").Value = 0.001771657596666667
$ws.Range("R2").Value = 0.01594491837
$ws.Range("S2").Value = 0.0162094769588191
$ws.Range("T2").Value = 0.0162094769588191

# Refresh the recomputed metrics for the new row 3 (FAPs -> FAPs)
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.9349228167457665
$ws.Range("P3").Value = 0.9349228167457664
$ws.Range("S3").Value = 0.9349228167457665
$ws.Range("T3").Value = 0.9349228167457664

# Refresh the recomputed metrics for the new row 4 (FAPs -> MuSCs)
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 0.2188686666666667
$ws.Range("N4").Value = 0.656606
$ws.Range("O4").Value = 0.04886770629541442
$ws.Range("P4").Value = 0.04886770629541441
$ws.Range("Q4").Value = 0.005341125028888889
$ws.Range("R4").Value = 0.04807012526
$ws.Range("S4").Value = 0.04886770629541442
$ws.Range("T4").Value = 0.04886770629541441
